{"js": "// Update the \"Perioadele campaniei...\" sentence throughout the document:\n// move \"2022\" after \"din\" and add \"pentru\" before \"Constela\u021bia Cygnus\".\n//   Before: \"Perioadele campaniei din Constela\u021bia Cygnus 2022: ...\"\n//   After:  \"Perioadele campaniei din 2022 pentru Constela\u021bia Cygnus: ...\"\nconst oldText = \"Perioadele campaniei din Constela\u021bia Cygnus 2022: 10-19 august, 9-18 septembrie, 8-17 octombrie\";\nconst newText = \"Perioadele campaniei din 2022 pentru Constela\u021bia Cygnus: 10-19 august, 9-18 septembrie, 8-17 octombrie\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the \"Perioadele campaniei...\" sentence throughout the document:\n# move \"2022\" after \"din\" and add \"pentru\" before \"Constela\u021bia Cygnus\".\n#   Before: \"Perioadele campaniei din Constela\u021bia Cygnus 2022: ...\"\n#   After:  \"Perioadele campaniei din 2022 pentru Constela\u021bia Cygnus: ...\"\n\n$d = $word.ActiveDocument\n\n$oldText = \"Perioadele campaniei din Constela\u021bia Cygnus 2022: 10-19 august, 9-18 septembrie, 8-17 octombrie\"\n$newText = \"Perioadele campaniei din 2022 pentru Constela\u021bia Cygnus: 10-19 august, 9-18 septembrie, 8-17 octombrie\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.Text = $newText\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\n    $oldText,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $newText,\n    2\n)\n"}
